$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data
$ws.Range("A2").Value = "Serbian First League"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-11-14"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "13:00:00"
$ws.Range("D2").Value = "FK Loznica"
$ws.Range("E2").Value = "Fk Smederevo"
$ws.Range("F2").Value = 1.01
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.01
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.02
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.01
$ws.Range("O2").Value = 1.01
$ws.Range("P2").Value = 1.1
$ws.Range("Q2").Value = 1.01
$ws.Range("R2").Value = 1.09
$ws.Range("S2").Value = 1.02
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3 data
$ws.Range("A3").Value = "Danish 2nd Division"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-11-14"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "14:00:00"
$ws.Range("D3").Value = "HIK Hellerup"
$ws.Range("E3").Value = "Vendsyssel FF"
$ws.Range("F3").Value = 2.84
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.26
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.7
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 1.86
$ws.Range("R3").Value = 1.24
$ws.Range("S3").Value = 1.86
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
